$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'57535"
$ws.Range("E2").Value = "'2802"
$ws.Range("A3").Value = "'65307"
$ws.Range("E3").Value = "'2579"
$ws.Range("A5").Value = "'51129"
$ws.Range("A7").Value = "'45879"
$ws.Range("E7").Value = "'3625"
$ws.Range("A8").Value = "'11046"
$ws.Range("E8").Value = "'5860"
$ws.Range("A9").Value = "'13753"
$ws.Range("E9").Value = "'5596"
$ws.Range("A10").Value = "'19617"
$ws.Range("B10").Value = "'38809086"
$ws.Range("C10").Value = "Kouenᶻᵍˣ"
$ws.Range("E10").Value = "'5205"
$ws.Range("A11").Value = "'20078"
$ws.Range("B11").Value = "'4756174"
$ws.Range("C11").Value = "純希です"
$ws.Range("E11").Value = "'5180"
$ws.Range("A12").Value = "'31056"
$ws.Range("E12").Value = "'4596"
$ws.Range("A13").Value = "'53777"
$ws.Range("E13").Value = "'2965"
$ws.Range("A14").Value = "'67104"
$ws.Range("E14").Value = "'2546"
$ws.Range("A16").Value = "'15220"
$ws.Range("E16").Value = "'5478"
$ws.Range("A17").Value = "'15974"
$ws.Range("E17").Value = "'5421"
$ws.Range("A18").Value = "'16327"
$ws.Range("E18").Value = "'5398"
$ws.Range("A19").Value = "'16379"
$ws.Range("E19").Value = "'5395"
$ws.Range("A20").Value = "'21354"
$ws.Range("B20").Value = "'55769051"
$ws.Range("C20").Value = "㊥叮叮当."
$ws.Range("E20").Value = "'5101"
$ws.Range("A21").Value = "'21413"
$ws.Range("B21").Value = "'46289694"
$ws.Range("C21").Value = "㊥Vincent"
$ws.Range("E21").Value = "'5098"
$ws.Range("A22").Value = "'21640"
$ws.Range("B22").Value = "'54085771"
$ws.Range("C22").Value = "㊥Matthieu"
$ws.Range("E22").Value = "'5085"
$ws.Range("A23").Value = "'31277"
$ws.Range("E23").Value = "'4587"
$ws.Range("A24").Value = "'32879"
$ws.Range("B24").Value = "'56585361"
$ws.Range("C24").Value = "`"㊥ go策划我要ali`""
$ws.Range("E24").Value = "'4519"
$ws.Range("A25").Value = "'33167"
$ws.Range("E25").Value = "'4507"
$ws.Range("A26").Value = "'33064"
$ws.Range("B26").Value = "'58839983"
$ws.Range("C26").Value = "每逢佳节胖六斤"
$ws.Range("E26").Value = "'4511"
$ws.Range("A27").Value = "'38182"
$ws.Range("B27").Value = "'55860890"
$ws.Range("C27").Value = "㊥Ethan"
$ws.Range("E27").Value = "'4246"
$ws.Range("A28").Value = "'39143"
$ws.Range("B28").Value = "'58408326"
$ws.Range("C28").Value = "`"Killer Bee`""
$ws.Range("E28").Value = "'4192"
$ws.Range("A29").Value = "'41362"
$ws.Range("B29").Value = "'1304123"
$ws.Range("C29").Value = "Cccccccccccc"
$ws.Range("E29").Value = "'4071"
$ws.Range("A30").Value = "'5761"
$ws.Range("E30").Value = "'6369"
$ws.Range("A31").Value = "'8381"
$ws.Range("E31").Value = "'6075"
$ws.Range("A32").Value = "'9819"
$ws.Range("B32").Value = "'11582001"
$ws.Range("C32").Value = "iMinatoX4"
$ws.Range("E32").Value = "'5959"
$ws.Range("A33").Value = "'11891"
$ws.Range("E33").Value = "'5779"
$ws.Range("A34").Value = "'12659"
$ws.Range("B34").Value = "'56133764"
$ws.Range("C34").Value = "ustcarter"
$ws.Range("E34").Value = "'5702"
$ws.Range("A35").Value = "'12820"
$ws.Range("B35").Value = "'55317038"
$ws.Range("C35").Value = "necman12345"
$ws.Range("E35").Value = "'5686"
$ws.Range("A36").Value = "'17604"
$ws.Range("E36").Value = "'5316"
$ws.Range("A37").Value = "'19041"
$ws.Range("E37").Value = "'5239"
$ws.Range("A38").Value = "'29886"
$ws.Range("E38").Value = "'4648"
$ws.Range("A39").Value = "'30627"
$ws.Range("E39").Value = "'4614"
$ws.Range("A40").Value = "'32291"
$ws.Range("E40").Value = "'4543"
$ws.Range("A41").Value = "'33771"
$ws.Range("B41").Value = "'38893233"
$ws.Range("C41").Value = "`"快乐 二哈`""
$ws.Range("E41").Value = "'4480"
$ws.Range("A42").Value = "'33805"
$ws.Range("B42").Value = "'52997727"
$ws.Range("C42").Value = "larios"
$ws.Range("E42").Value = "'4478"
$ws.Range("A43").Value = "'33904"
$ws.Range("B43").Value = "'56379103"
$ws.Range("C43").Value = "Globalking"
$ws.Range("E43").Value = "'4472"
$ws.Range("A44").Value = "'35726"
$ws.Range("B44").Value = "'32316256"
$ws.Range("C44").Value = "`"秋の風 ..`""
$ws.Range("E44").Value = "'4377"
$ws.Range("A45").Value = "'36794"
$ws.Range("B45").Value = "'50837459"
$ws.Range("C45").Value = "NINE日"
$ws.Range("E45").Value = "'4316"
$ws.Range("A46").Value = "'39798"
$ws.Range("B46").Value = "'55634661"
$ws.Range("C46").Value = "Opalus"
$ws.Range("E46").Value = "'4158"
$ws.Range("A47").Value = "'41067"
$ws.Range("B47").Value = "'58203298"
$ws.Range("C47").Value = "权旨qua"
$ws.Range("E47").Value = "'4087"
$ws.Range("A48").Value = "'41580"
$ws.Range("B48").Value = "'59020292"
$ws.Range("C48").Value = "Sharnoth"
$ws.Range("E48").Value = "'4056"
$ws.Range("A49").Value = "'48765"
$ws.Range("E49").Value = "'3322"
$ws.Range("A50").Value = "'57011"
$ws.Range("E50").Value = "'2821"
$ws.Range("A51").Value = "'67993"
$ws.Range("A52").Value = "'61840"
$ws.Range("E52").Value = "'2661"
$ws.Range("A53").Value = "'50456"
$ws.Range("E53").Value = "'3181"
$ws.Range("A56").Value = "'43093"
$ws.Range("E56").Value = "'3982"
$ws.Range("A57").Value = "'47058"
$ws.Range("E57").Value = "'3497"
$ws.Range("A58").Value = "'59594"
$ws.Range("E58").Value = "'2729"
$ws.Range("A59").Value = "'68751"
$ws.Range("E59").Value = "'2520"
$ws.Range("A60").Value = "'108824"
$ws.Range("B60").Value = "'15436348"
$ws.Range("C60").Value = "Lucas"
$ws.Range("E60").Value = "'1498"
$ws.Range("A61").Value = "'108977"
$ws.Range("B61").Value = "'57219176"
$ws.Range("C61").Value = "青莲道人"
$ws.Range("E61").Value = "'1497"
$ws.Range("A62").Value = "'111104"
$ws.Range("A63").Value = "'124017"
$ws.Range("E63").Value = "'1300"
$ws.Range("A75").Value = "'46855"
$ws.Range("E75").Value = "'3519"
$ws.Range("A78").Value = "'90397"
$ws.Range("E78").Value = "'1911"
$ws.Range("A79").Value = "'97530"
$ws.Range("A80").Value = "'162549"
$ws.Range("A81").Value = "'217788"
